# Fruta / hortaliza, semanal
#
# The raw daily price records (rows 2-20) are re-ordered/updated to reflect
# the weekly aggregation. For each destination row, the Fecha (D), Volumen
# (M), Precio minimo (N), Precio maximo (O), Precio promedio ponderado (P)
# and Precio $/Kg (S) values are replaced with the values taken from the
# row indicated below (all other columns are left untouched).
#
# destination row -> source row (values copied from source's original data)
#   2 <- 7      8 <- 20     14 <- 9      20 <- 11
#   3 <- 2      9 <- 5      15 <- 14
#   4 <- 3     10 <- 16     16 <- 18
#   5 <- 19    11 <- 4      17 <- 13
#   6 <- 6     12 <- 15     18 <- 8
#   7 <- 17    13 <- 12     19 <- 10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original values of the columns that move, keyed by row number,
# before any of them get overwritten. Value2 is used (rather than Value) when
# reading/round-tripping through a variable so the raw number is preserved.
$original = @{}
for ($row = 2; $row -le 20; $row++) {
    $original[$row] = @{
        D = $ws.Range("D$row").Value2
        M = $ws.Range("M$row").Value2
        N = $ws.Range("N$row").Value2
        O = $ws.Range("O$row").Value2
        P = $ws.Range("P$row").Value2
        S = $ws.Range("S$row").Value2
    }
}

$rowMap = @{
    2 = 7;  3 = 2;  4 = 3;  5 = 19; 6 = 6;  7 = 17; 8 = 20; 9 = 5;  10 = 16
    11 = 4; 12 = 15; 13 = 12; 14 = 9; 15 = 14; 16 = 18; 17 = 13; 18 = 8
    19 = 10; 20 = 11
}

foreach ($destRow in $rowMap.Keys | Sort-Object) {
    $srcRow = $rowMap[$destRow]
    $src = $original[$srcRow]

    $ws.Range("D$destRow").Value2 = $src.D
    $ws.Range("M$destRow").Value2 = $src.M
    $ws.Range("N$destRow").Value2 = $src.N
    $ws.Range("O$destRow").Value2 = $src.O
    $ws.Range("P$destRow").Value2 = $src.P
    $ws.Range("S$destRow").Value2 = $src.S
}
